$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7342.5
$ws.Range("I64").Value = 5985
$ws.Range("J64").Value = 7924.2856
$ws.Range("K64").Value = 5985
$ws.Range("L64").Value = 7924.2856
$ws.Range("M64").Value = -5737
$ws.Range("N64").Value = -8420.285599999999
$ws.Range("H67").Value = 7342.5
$ws.Range("I67").Value = 5985
$ws.Range("J67").Value = 7924.2856
$ws.Range("K67").Value = 5985
$ws.Range("L67").Value = 7924.2856
$ws.Range("M67").Value = -5127
$ws.Range("N67").Value = -9640.285599999999
$ws.Range("H74").Value = 7863.7
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 7863.7
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 7863.7
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -9735.700000000001
$ws.Range("H77").Value = 7863.7
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 7863.7
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 39318.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -48678.5
$ws.Range("H132").Value = 43481770
$ws.Range("I132").Value = 45458172
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 136374516
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -136371986
$ws.Range("N132").Value = -7760
$ws.Range("H137").Value = 83652.05
$ws.Range("I137").Value = 121078.07
$ws.Range("J137").Value = 3453.4285
$ws.Range("K137").Value = 363234.21
$ws.Range("L137").Value = 10360.2855
$ws.Range("M137").Value = -360684.21
$ws.Range("N137").Value = -15460.2855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7105.041
$ws.Range("I32").Value = 3631.9333
$ws.Range("K32").Value = 3631.9333
$ws.Range("M32").Value = -3344.9333
$ws.Range("H39").Value = 11449.5
$ws.Range("J39").Value = 11900
$ws.Range("L39").Value = 11900
$ws.Range("N39").Value = -12940
$ws.Range("H119").Value = 56899.91
$ws.Range("J119").Value = 56899.91
$ws.Range("L119").Value = 56899.91
$ws.Range("N119").Value = -66575.91
$ws.Range("H122").Value = 17366652
$ws.Range("I122").Value = 55558204
$ws.Range("J122").Value = 2090029.8
$ws.Range("K122").Value = 166674612
$ws.Range("L122").Value = 6270089.4
$ws.Range("M122").Value = -166672162
$ws.Range("N122").Value = -6274989.4
$ws.Range("H132").Value = 2622.6487
$ws.Range("I132").Value = 1858.3871
$ws.Range("J132").Value = 6571.3335
$ws.Range("K132").Value = 5575.1613
$ws.Range("L132").Value = 19714.0005
$ws.Range("M132").Value = -3045.1613
$ws.Range("N132").Value = -24774.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 40401
$ws.Range("J92").Value = 40401
$ws.Range("L92").Value = 40401
$ws.Range("N92").Value = -45393

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1420.4
$ws.Range("I105").Value = 1034
$ws.Range("K105").Value = 1034
$ws.Range("M105").Value = 713
$ws.Range("H112").Value = 26999
$ws.Range("J112").Value = 26999
$ws.Range("L112").Value = 26999
$ws.Range("N112").Value = -29953
$ws.Range("H133").Value = 69949.25
$ws.Range("J133").Value = 69949.25
$ws.Range("L133").Value = 69949.25
$ws.Range("N133").Value = -75009.25
$ws.Range("H141").Value = 204448.45
$ws.Range("J141").Value = 204448.45
$ws.Range("L141").Value = 204448.45
$ws.Range("N141").Value = -214808.45

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 59854.65
$ws.Range("J5").Value = 126074.5
$ws.Range("L5").Value = 378223.5
$ws.Range("N5").Value = -378447.5
$ws.Range("H12").Value = 81988.27
$ws.Range("I12").Value = 177936.4
$ws.Range("J12").Value = 2031.5
$ws.Range("K12").Value = 533809.2
$ws.Range("L12").Value = 6094.5
$ws.Range("M12").Value = -533636.2
$ws.Range("N12").Value = -6440.5
$ws.Range("H68").Value = 1549.9286
$ws.Range("I68").Value = 701
$ws.Range("J68").Value = 2186.625
$ws.Range("K68").Value = 2103
$ws.Range("L68").Value = 6559.875
$ws.Range("M68").Value = -1292
$ws.Range("N68").Value = -8181.875
$ws.Range("H71").Value = 1549.9286
$ws.Range("I71").Value = 701
$ws.Range("J71").Value = 2186.625
$ws.Range("K71").Value = 6309
$ws.Range("L71").Value = 19679.625
$ws.Range("M71").Value = -2253
$ws.Range("N71").Value = -27791.625
$ws.Range("H107").Value = 1107.9445
$ws.Range("I107").Value = 1022.6923
$ws.Range("K107").Value = 3068.0769
$ws.Range("M107").Value = -1148.0769
$ws.Range("H130").Value = 2847.6
$ws.Range("I130").Value = 2568
$ws.Range("K130").Value = 7704
$ws.Range("M130").Value = -2684
$ws.Range("H131").Value = 10420215
$ws.Range("J131").Value = 11498270
$ws.Range("L131").Value = 34494810
$ws.Range("N131").Value = -34504890
$ws.Range("H132").Value = 1678.909
$ws.Range("I132").Value = 1499
$ws.Range("K132").Value = 13491
$ws.Range("M132").Value = -10961
$ws.Range("H135").Value = 59854.65
$ws.Range("J135").Value = 126074.5
$ws.Range("L135").Value = 1134670.5
$ws.Range("N135").Value = -1139740.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4599
$ws.Range("J11").Value = 4599
$ws.Range("L11").Value = 4599
$ws.Range("N11").Value = -4877
$ws.Range("H42").Value = 30972
$ws.Range("J42").Value = 30972
$ws.Range("L42").Value = 30972
$ws.Range("N42").Value = -31942
$ws.Range("H115").Value = 30972
$ws.Range("J115").Value = 30972
$ws.Range("L115").Value = 30972
$ws.Range("N115").Value = -33322
$ws.Range("H122").Value = 470339.53
$ws.Range("I122").Value = 594636.4399999999
$ws.Range("J122").Value = 4226
$ws.Range("K122").Value = 1783909.32
$ws.Range("L122").Value = 12678
$ws.Range("M122").Value = -1781459.32
$ws.Range("N122").Value = -17578
$ws.Range("H136").Value = 15680
$ws.Range("J136").Value = 16702.5
$ws.Range("L136").Value = 50107.5
$ws.Range("N136").Value = -55207.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4809.087
$ws.Range("I7").Value = 3410.5334
$ws.Range("K7").Value = 3410.5334
$ws.Range("M7").Value = -3298.5334
$ws.Range("H16").Value = 705.1111
$ws.Range("I16").Value = 612.5
$ws.Range("J16").Value = 890.3333
$ws.Range("K16").Value = 612.5
$ws.Range("L16").Value = 890.3333
$ws.Range("M16").Value = -442.5
$ws.Range("N16").Value = -1230.3333
$ws.Range("H32").Value = 11285.25
$ws.Range("I32").Value = 11285.25
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11285.25
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -10968.25
$ws.Range("N32").ClearContents()
$ws.Range("H46").Value = 7198.75
$ws.Range("I46").Value = 2597.25
$ws.Range("J46").Value = 9499.5
$ws.Range("K46").Value = 2597.25
$ws.Range("L46").Value = 9499.5
$ws.Range("M46").Value = -2409.25
$ws.Range("N46").Value = -9875.5
$ws.Range("H61").Value = 12350213
$ws.Range("I61").Value = 15874273
$ws.Range("K61").Value = 15874273
$ws.Range("M61").Value = -15874071
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256
$ws.Range("H113").Value = 12350213
$ws.Range("I113").Value = 15874273
$ws.Range("K113").Value = 15874273
$ws.Range("M113").Value = -15872103
$ws.Range("H126").Value = 4809.087
$ws.Range("I126").Value = 3410.5334
$ws.Range("K126").Value = 10231.6002
$ws.Range("M126").Value = -7761.600199999999
$ws.Range("H132").Value = 8228.538
$ws.Range("I132").Value = 9385.223
$ws.Range("K132").Value = 28155.669
$ws.Range("M132").Value = -25625.669
$ws.Range("H136").Value = 49050.867
$ws.Range("I136").Value = 114148.72
$ws.Range("J136").Value = 5652.2964
$ws.Range("K136").Value = 342446.16
$ws.Range("L136").Value = 16956.8892
$ws.Range("M136").Value = -339896.16
$ws.Range("N136").Value = -22056.8892

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 16500
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = -2750
$ws.Range("N33").Value = -30500
$ws.Range("H36").Value = 16500
$ws.Range("I36").Value = 3000
$ws.Range("J36").Value = 30000
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 30000
$ws.Range("M36").Value = -2750
$ws.Range("N36").Value = -30500
$ws.Range("H94").Value = 24996
$ws.Range("J94").Value = 24996
$ws.Range("L94").Value = 24996
$ws.Range("N94").Value = -26798
$ws.Range("H126").Value = 1286.5
$ws.Range("I126").Value = 1286.8334
$ws.Range("K126").Value = 3860.5002
$ws.Range("M126").Value = -1390.5002
$ws.Range("H132").Value = 16852060
$ws.Range("I132").Value = 19611956
$ws.Range("K132").Value = 58835868
$ws.Range("M132").Value = -58833338
$ws.Range("H136").Value = 3918.0605
$ws.Range("I136").Value = 3577.077
$ws.Range("K136").Value = 10731.231
$ws.Range("M136").Value = -8181.231
